# Auto-generated edit script: updates odds data for Jogos da Semana FlashScore sheet
# Applies the value changes captured in the commit diff (various odds/columns updated).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 2.45
$ws.Range("I3").Value = 3.4
$ws.Range("L3").Value = 4.33
$ws.Range("M3").Value = 1.11
$ws.Range("O3").Value = 1.54
$ws.Range("X3").Value = 1.1
$ws.Range("AQ3").Value = 41
# Row 4
$ws.Range("M4").Value = 1.13
$ws.Range("O4").Value = 1.63
$ws.Range("X4").Value = 1.07
$ws.Range("Z4").Value = 2.08
# Row 7
$ws.Range("G7").Value = 1.72
$ws.Range("H7").Value = 3.15
$ws.Range("J7").Value = 2.27
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 5.8
$ws.Range("M7").Value = 1.11
$ws.Range("N7").Value = 5.7
$ws.Range("P7").Value = 2.55
$ws.Range("S7").Value = 2.35
$ws.Range("T7").Value = 1.53
$ws.Range("W7").Value = 4.1
$ws.Range("X7").Value = 1.19
$ws.Range("Y7").Value = 1.5
$ws.Range("Z7").Value = 2.42
$ws.Range("AC7").Value = 5.1
$ws.Range("AI7").Value = 5.7
# Row 8
$ws.Range("H8").Value = 2.55
$ws.Range("J8").Value = 3.75
$ws.Range("L8").Value = 3.7
$ws.Range("N8").Value = 4.3
$ws.Range("O8").Value = 1.72
$ws.Range("P8").Value = 2
$ws.Range("S8").Value = 3.1
$ws.Range("T8").Value = 1.32
$ws.Range("W8").Value = 5.8
$ws.Range("AA8").Value = 2.35
$ws.Range("AB8").Value = 1.53
$ws.Range("AD8").Value = 13
$ws.Range("AE8").Value = 12
$ws.Range("AG8").Value = 37
$ws.Range("AH8").Value = 65
$ws.Range("AI8").Value = 4.3
$ws.Range("AN8").Value = 5.7
$ws.Range("AO8").Value = 12.5
# Row 9
$ws.Range("H9").Value = 4.1
$ws.Range("I9").Value = 8.75
$ws.Range("J9").Value = 1.93
$ws.Range("K9").Value = 2.18
$ws.Range("L9").Value = 7.5
$ws.Range("M9").Value = 1.06
$ws.Range("N9").Value = 7.3
$ws.Range("O9").Value = 1.29
$ws.Range("P9").Value = 3.3
$ws.Range("S9").Value = 1.85
$ws.Range("T9").Value = 1.85
$ws.Range("W9").Value = 3.05
$ws.Range("X9").Value = 1.33
$ws.Range("Y9").Value = 1.42
$ws.Range("Z9").Value = 2.67
$ws.Range("AA9").Value = 2.12
$ws.Range("AB9").Value = 1.65
$ws.Range("AC9").Value = 5.8
$ws.Range("AD9").Value = 5.9
$ws.Range("AE9").Value = 8.25
$ws.Range("AF9").Value = 8.5
$ws.Range("AG9").Value = 12
$ws.Range("AH9").Value = 32
$ws.Range("AI9").Value = 7.3
$ws.Range("AK9").Value = 21
$ws.Range("AL9").Value = 110
$ws.Range("AM9").Value = 900
$ws.Range("AN9").Value = 20
$ws.Range("AO9").Value = 60
$ws.Range("AP9").Value = 26
$ws.Range("AR9").Value = 110
$ws.Range("AS9").Value = 90
# Row 10
$ws.Range("G10").Value = 1.85
$ws.Range("H10").Value = 3.7
$ws.Range("S10").Value = 1.5
$ws.Range("T10").Value = 2.5
$ws.Range("AA10").Value = 1.5
$ws.Range("AB10").Value = 2.5
$ws.Range("AD10").Value = 12
$ws.Range("AF10").Value = 19
$ws.Range("AH10").Value = 19
$ws.Range("AI10").Value = 19
$ws.Range("AJ10").Value = 8
$ws.Range("AR10").Value = 23
# Row 12
$ws.Range("G12").Value = 1.18
$ws.Range("H12").Value = 6.5
$ws.Range("I12").Value = 11
$ws.Range("J12").Value = 1.53
$ws.Range("L12").Value = 9.5
$ws.Range("M12").Value = 1.01
$ws.Range("N12").Value = 17
$ws.Range("AA12").Value = 2.2
$ws.Range("AB12").Value = 1.62
$ws.Range("AE12").Value = 11
$ws.Range("AF12").Value = 7
$ws.Range("AK12").Value = 29
$ws.Range("AL12").Value = 81
$ws.Range("AP12").Value = 34
$ws.Range("AQ12").Value = 151
$ws.Range("AR12").Value = 81
$ws.Range("AS12").Value = 67
# Row 13
$ws.Range("S13").Value = 1.53
$ws.Range("T13").Value = 2.38
$ws.Range("U13").Value = 1.85
$ws.Range("V13").Value = 1.95
$ws.Range("W13").Value = 2.25
$ws.Range("X13").Value = 1.57
# Row 14
$ws.Range("S14").Value = 1.8
$ws.Range("T14").Value = 2
$ws.Range("AM14").Value = 800
# Row 15
$ws.Range("G15").Value = 1.73
$ws.Range("H15").Value = 4.1
$ws.Range("I15").Value = 4.2
$ws.Range("J15").Value = 2.25
$ws.Range("L15").Value = 4.33
$ws.Range("M15").Value = 1.03
$ws.Range("N15").Value = 17
$ws.Range("U15").Value = 1.8
$ws.Range("V15").Value = 2.05
$ws.Range("W15").Value = 2.2
$ws.Range("X15").Value = 1.62
$ws.Range("AA15").Value = 1.53
$ws.Range("AB15").Value = 2.38
$ws.Range("AD15").Value = 10
$ws.Range("AF15").Value = 15
$ws.Range("AH15").Value = 19
$ws.Range("AJ15").Value = 8
$ws.Range("AK15").Value = 12
$ws.Range("AL15").Value = 34
$ws.Range("AM15").Value = 101
$ws.Range("AQ15").Value = 41
# Row 16
$ws.Range("G16").Value = 3.4
$ws.Range("I16").Value = 2.1
$ws.Range("J16").Value = 3.75
$ws.Range("L16").Value = 2.75
$ws.Range("N16").Value = 12
$ws.Range("AC16").Value = 12
$ws.Range("AF16").Value = 34
$ws.Range("AH16").Value = 29
$ws.Range("AI16").Value = 12
$ws.Range("AM16").Value = 151
$ws.Range("AO16").Value = 11
$ws.Range("AR16").Value = 17
# Row 17
$ws.Range("G17").Value = 5.2
$ws.Range("H17").Value = 4.15
$ws.Range("J17").Value = 5
$ws.Range("K17").Value = 2.35
$ws.Range("N17").Value = 8.5
$ws.Range("P17").Value = 3.9
$ws.Range("W17").Value = 2.57
$ws.Range("X17").Value = 1.45
$ws.Range("AC17").Value = 16
$ws.Range("AF17").Value = 90
$ws.Range("AG17").Value = 45
$ws.Range("AI17").Value = 8.5
$ws.Range("AJ17").Value = 8
$ws.Range("AK17").Value = 16
$ws.Range("AN17").Value = 7.8
# Row 18
$ws.Range("G18").Value = 1.65
$ws.Range("H18").Value = 3.8
$ws.Range("I18").Value = 4.5
$ws.Range("J18").Value = 2.22
$ws.Range("K18").Value = 2.2
$ws.Range("L18").Value = 4.7
$ws.Range("O18").Value = 1.27
$ws.Range("P18").Value = 3.1
$ws.Range("S18").Value = 1.8
$ws.Range("T18").Value = 1.8
$ws.Range("W18").Value = 2.87
$ws.Range("X18").Value = 1.31
$ws.Range("AA18").Value = 1.82
$ws.Range("AB18").Value = 1.8
$ws.Range("AC18").Value = 6.8
$ws.Range("AD18").Value = 7.6
$ws.Range("AE18").Value = 8.25
$ws.Range("AF18").Value = 12.5
$ws.Range("AG18").Value = 13.5
$ws.Range("AI18").Value = 10.5
$ws.Range("AJ18").Value = 7.4
$ws.Range("AK18").Value = 17
$ws.Range("AL18").Value = 80
$ws.Range("AM18").Value = 700
$ws.Range("AN18").Value = 12.5
$ws.Range("AO18").Value = 25
$ws.Range("AP18").Value = 15
$ws.Range("AQ18").Value = 75
$ws.Range("AR18").Value = 45
$ws.Range("AS18").Value = 50
